{"js": "// Replace the date line and the 25 multiplication problems with their\n// updated values, per the commit diff. Every \"old\" string below occurs\n// exactly once in the document, so a body.search + insertText(\"Replace\")\n// round-trip is safe and keeps all run/paragraph formatting untouched.\nconst replacements = [\n  [\"2023-11-20 Monday\", \"2023-11-21 Tuesday\"],\n  [\"69\u00d726=\", \"57\u00d754=\"],\n  [\"98\u00d768=\", \"28\u00d733=\"],\n  [\"23\u00d751=\", \"46\u00d784=\"],\n  [\"36\u00d715=\", \"42\u00d716=\"],\n  [\"26\u00d775=\", \"17\u00d739=\"],\n  [\"97\u00d733=\", \"72\u00d720=\"],\n  [\"75\u00d731=\", \"51\u00d743=\"],\n  [\"92\u00d713=\", \"19\u00d719=\"],\n  [\"73\u00d768=\", \"53\u00d731=\"],\n  [\"16\u00d753=\", \"56\u00d798=\"],\n  [\"72\u00d742=\", \"89\u00d760=\"],\n  [\"64\u00d750=\", \"31\u00d718=\"],\n  [\"70\u00d781=\", \"67\u00d721=\"],\n  [\"45\u00d738=\", \"63\u00d789=\"],\n  [\"27\u00d797=\", \"52\u00d741=\"],\n  [\"46\u00d737=\", \"48\u00d790=\"],\n  [\"47\u00d756=\", \"40\u00d737=\"],\n  [\"19\u00d715=\", \"91\u00d756=\"],\n  [\"18\u00d750=\", \"81\u00d722=\"],\n  [\"38\u00d745=\", \"49\u00d765=\"],\n  [\"50\u00d771=\", \"96\u00d726=\"],\n  [\"16\u00d785=\", \"99\u00d739=\"],\n  [\"49\u00d740=\", \"30\u00d788=\"],\n  [\"41\u00d725=\", \"94\u00d722=\"],\n  [\"27\u00d737=\", \"48\u00d761=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems in place.\n# Each \"old\" string occurs exactly once in the document, so Find/Replace\n# with wdReplaceAll is a safe 1-for-1 substitution that keeps all run /\n# paragraph formatting untouched.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"2023-11-20 Monday\", $false, $false, $false, $false, $false, $true, 1, $false, \"2023-11-21 Tuesday\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"69\u00d726=\", $false, $false, $false, $false, $false, $true, 1, $false, \"57\u00d754=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"98\u00d768=\", $false, $false, $false, $false, $false, $true, 1, $false, \"28\u00d733=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"23\u00d751=\", $false, $false, $false, $false, $false, $true, 1, $false, \"46\u00d784=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"36\u00d715=\", $false, $false, $false, $false, $false, $true, 1, $false, \"42\u00d716=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"26\u00d775=\", $false, $false, $false, $false, $false, $true, 1, $false, \"17\u00d739=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"97\u00d733=\", $false, $false, $false, $false, $false, $true, 1, $false, \"72\u00d720=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"75\u00d731=\", $false, $false, $false, $false, $false, $true, 1, $false, \"51\u00d743=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"92\u00d713=\", $false, $false, $false, $false, $false, $true, 1, $false, \"19\u00d719=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"73\u00d768=\", $false, $false, $false, $false, $false, $true, 1, $false, \"53\u00d731=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"16\u00d753=\", $false, $false, $false, $false, $false, $true, 1, $false, \"56\u00d798=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"72\u00d742=\", $false, $false, $false, $false, $false, $true, 1, $false, \"89\u00d760=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"64\u00d750=\", $false, $false, $false, $false, $false, $true, 1, $false, \"31\u00d718=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"70\u00d781=\", $false, $false, $false, $false, $false, $true, 1, $false, \"67\u00d721=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"45\u00d738=\", $false, $false, $false, $false, $false, $true, 1, $false, \"63\u00d789=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"27\u00d797=\", $false, $false, $false, $false, $false, $true, 1, $false, \"52\u00d741=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"46\u00d737=\", $false, $false, $false, $false, $false, $true, 1, $false, \"48\u00d790=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"47\u00d756=\", $false, $false, $false, $false, $false, $true, 1, $false, \"40\u00d737=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"19\u00d715=\", $false, $false, $false, $false, $false, $true, 1, $false, \"91\u00d756=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"18\u00d750=\", $false, $false, $false, $false, $false, $true, 1, $false, \"81\u00d722=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"38\u00d745=\", $false, $false, $false, $false, $false, $true, 1, $false, \"49\u00d765=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"50\u00d771=\", $false, $false, $false, $false, $false, $true, 1, $false, \"96\u00d726=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"16\u00d785=\", $false, $false, $false, $false, $false, $true, 1, $false, \"99\u00d739=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"49\u00d740=\", $false, $false, $false, $false, $false, $true, 1, $false, \"30\u00d788=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"41\u00d725=\", $false, $false, $false, $false, $false, $true, 1, $false, \"94\u00d722=\", 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"27\u00d737=\", $false, $false, $false, $false, $false, $true, 1, $false, \"48\u00d761=\", 2)\n\n"}
